$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on columns D (Price) and G (Hora) so the numeric-looking
# strings are stored as text (matching the source data feed format) instead of
# being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Updated Price (column D) values
$ws.Cells.Item(2, 4).Value = "244.80"
$ws.Cells.Item(3, 4).Value = "23.13"
$ws.Cells.Item(4, 4).Value = "5.441"
$ws.Cells.Item(5, 4).Value = "0.05972"
$ws.Cells.Item(6, 4).Value = "3.390"
$ws.Cells.Item(8, 4).Value = "0.9269"
$ws.Cells.Item(9, 4).Value = "0.1430"
$ws.Cells.Item(11, 4).Value = "0.03374"
$ws.Cells.Item(13, 4).Value = "0.09352"
$ws.Cells.Item(14, 4).Value = "3.936"
$ws.Cells.Item(15, 4).Value = "0.001606"
$ws.Cells.Item(16, 4).Value = "0.04825"
$ws.Cells.Item(17, 4).Value = "0.0005943"
$ws.Cells.Item(18, 4).Value = "0.005601"
$ws.Cells.Item(19, 4).Value = "0.004157"
$ws.Cells.Item(20, 4).Value = "0.0009838"
$ws.Cells.Item(23, 4).Value = "6.459"
$ws.Cells.Item(26, 4).Value = "0.1341"
$ws.Cells.Item(27, 4).Value = "0.0002447"
$ws.Cells.Item(40, 4).Value = "0.03940"
$ws.Cells.Item(41, 4).Value = "0.1074"
$ws.Cells.Item(42, 4).Value = "0.002682"
$ws.Cells.Item(43, 4).Value = "0.006184"
$ws.Cells.Item(44, 4).Value = "0.007241"
$ws.Cells.Item(45, 4).Value = "0.00005125"
$ws.Cells.Item(47, 4).Value = "0.0005803"
$ws.Cells.Item(49, 4).Value = "0.002277"

# Updated Hora (column G) values: all rows bump from "10" to "11"
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 7).Value = "11"
}
